$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.439.00"
$ws.Range("D2").ClearFormats()

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.092.10"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.10%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.99"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.24%  "

$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5214"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.27%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4374"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.24"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +16.83%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08849"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.152"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.58%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.24"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.83%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.087.11"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.700"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.689"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "95.80"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001120"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06593"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.86%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.23"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.19%  "

$ws.Range("E21").Value = "  +0.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.272"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.480.08"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.26"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.43%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.340"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.69%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.318.85"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.25"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.567"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.29"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "131.67"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.59%  "

$ws.Range("E31").Value = "  -0.58%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1067"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.646"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +7.64%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.162"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.51%  "

$ws.Range("E35").Value = "  -0.34%  "

$ws.Range("E36").Value = "  +4.46%  "

$ws.Range("E37").Value = "  -0.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06804"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.461"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.33%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.64"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2255"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.46%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6879"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.256"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("E44").Value = "  +0.19%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6347"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.92"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.19%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.193"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.99%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.624"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.39%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.237"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +7.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.245"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.79"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.83%  "
